# Upgrade: -Cambio Index da vin a targa; -Cambio numero minimo di posti
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the INDEX key from "numero_telaio" (VIN) row to "targa" (license plate) row,
# and change numero_telaio's attribute from "NOT NULL,UNIQUE" to "NULL".
$ws.Range("E10").Value = "INDEX"
$ws.Range("E10").HorizontalAlignment = -4108
$ws.Range("D11").Value = "NULL"
$ws.Range("E11").Value = ""

# Change the minimum default number of occupants from 1 to 2.
$ws.Range("D18").Value = "NOT NULL, DEFAULT (2)"
